$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.349.04'
$ws.Range("E2").Value = '  -2.28%  '

$ws.Range("D3").Value = '2.644.27'

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '598.74'
$ws.Range("E5").Value = '  -0.96%  '

$ws.Range("D6").Value = '165.82'
$ws.Range("E6").Value = '  -1.98%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  -0.95%  '

$ws.Range("D9").Value = '2.643.33'
$ws.Range("E9").Value = '  -3.42%  '

$ws.Range("E10").Value = '  -0.81%  '

$ws.Range("E11").Value = '  +1.40%  '

$ws.Range("E12").Value = '  -0.73%  '

$ws.Range("D13").Value = '5.22'
$ws.Range("E13").Value = '  -2.17%  '

$ws.Range("D14").Value = '27.97'
$ws.Range("E14").Value = '  -2.71%  '

$ws.Range("D15").Value = '3.123.41'
$ws.Range("E15").Value = '  -3.41%  '

$ws.Range("E16").Value = '  -3.21%  '

$ws.Range("D17").Value = '67.293.54'
$ws.Range("E17").Value = '  -2.06%  '

$ws.Range("D18").Value = '2.631.49'
$ws.Range("E18").Value = '  -4.21%  '

$ws.Range("D19").Value = '11.89'
$ws.Range("E19").Value = '  +0.40%  '

$ws.Range("D20").Value = '7.87'
$ws.Range("E20").Value = '  +2.03%  '

$ws.Range("D21").Value = '363.72'
$ws.Range("E21").Value = '  -2.87%  '

$ws.Range("E22").Value = '  -3.05%  '

$ws.Range("E23").Value = '  -3.43%  '

$ws.Range("D24").Value = '11.08'
$ws.Range("E24").Value = '  +11.03%  '

$ws.Range("E25").Value = '  -6.68%  '

$ws.Range("E26").Value = '  +0.06%  '

$ws.Range("D27").Value = '70.85'
$ws.Range("E27").Value = '  -4.11%  '

$ws.Range("D28").Value = '2.778.17'

$ws.Range("E29").Value = '  -4.11%  '

$ws.Range("E30").Value = '  +0.26%  '

$ws.Range("D31").Value = '554.12'
$ws.Range("E31").Value = '  -6.88%  '

$ws.Range("E32").Value = '  -3.33%  '

$ws.Range("E33").Value = '  -4.22%  '

$ws.Range("D34").Value = '1.93'
$ws.Range("E34").Value = '  -1.99%  '

$ws.Range("E35").Value = '  -0.44%  '

$ws.Range("E36").Value = '  -0.06%  '

$ws.Range("E37").Value = '  -5.34%  '

$ws.Range("D38").Value = '157.64'
$ws.Range("E38").Value = '  -2.39%  '

$ws.Range("D39").Value = '19.42'
$ws.Range("E39").Value = '  -2.25%  '

$ws.Range("E40").Value = '  -2.49%  '

$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '1.82'
$ws.Range("E41").Value = '  -5.30%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D42").Value = '5.26'
$ws.Range("E42").Value = '  -4.52%  '

$ws.Range("D43").Value = '17.91'
$ws.Range("E43").Value = '  -0.48%  '

$ws.Range("E44").Value = '  -6.18%  '

$ws.Range("E45").Value = '  -0.01%  '

$ws.Range("D46").Value = '40.16'
$ws.Range("E46").Value = '  -1.87%  '

$ws.Range("E47").Value = '  -3.81%  '

$ws.Range("E48").Value = '  -1.63%  '

$ws.Range("D49").Value = '154.23'
$ws.Range("E49").Value = '  -1.77%  '

$ws.Range("E51").Value = '  -3.77%  '
